$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The shared string "2016-10-20 01:43:30" (Latest HO Xliff Generate Date /
# Latest Handoff Datetime) is reused by several rows that happen to share
# the same timestamp; updating it updates all of its occurrences:
#   Overview!G4:G7  and  de-de!H4:H7
$wsOverview.Range("G4").Value = "2016-10-20 01:44:16"
$wsOverview.Range("G5").Value = "2016-10-20 01:44:16"
$wsOverview.Range("G6").Value = "2016-10-20 01:44:16"
$wsOverview.Range("G7").Value = "2016-10-20 01:44:16"

$wsDeDe.Range("H4").Value = "2016-10-20 01:44:16"
$wsDeDe.Range("H5").Value = "2016-10-20 01:44:16"
$wsDeDe.Range("H6").Value = "2016-10-20 01:44:16"
$wsDeDe.Range("H7").Value = "2016-10-20 01:44:16"

# zh-cn table: Priority low -> ht for rows 4-7
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

# zh-cn table: Latest Handoff Datetime update for rows 4-7
# ("2016-10-20 01:43:18" -> "2016-10-20 01:44:04", only used here)
$wsZhCn.Range("H4").Value = "2016-10-20 01:44:04"
$wsZhCn.Range("H5").Value = "2016-10-20 01:44:04"
$wsZhCn.Range("H6").Value = "2016-10-20 01:44:04"
$wsZhCn.Range("H7").Value = "2016-10-20 01:44:04"

# de-de table: Priority low -> ht for rows 4-7
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"
